# ejercicios en linea lab 1 listos
#
# The first heading is renumbered from "2." to "3." - i.e. the author
# selected the "2" in "          2. RECURSIVE DIGIT SUM" and typed "3"
# over it. That kind of edit leaves Word's auto-managed "_GoBack"
# bookmark sitting right after the freshly typed character, which also
# relocates the bookmark away from wherever it used to be (Word keeps
# only one "_GoBack" bookmark at a time), and it leaves the newly typed
# text in its own run, distinct from the run(s) that used to surround it.

$d = $word.ActiveDocument

# Locate the "2." that starts the first heading ("          2. RECURSIVE
# DIGIT SUM"); Find.Execute collapses/extends $rng onto the match so we
# can read exact character offsets back out of it instead of hard-coding
# paragraph-relative numbers.
$rng = $d.Content
$rng.Find.Execute("2.") | Out-Null

$digitStart = $rng.Start
$digitEnd = $digitStart + 1

# Mark the boundary just before the digit with a throwaway bookmark so
# the leading spaces stay in their own run, separate from the digit that
# is about to be retyped.
$preSplit = $d.Range($digitStart, $digitStart)
$d.Bookmarks.Add("zzTempSplit", $preSplit) | Out-Null

# Drop the "_GoBack" bookmark right after the digit, where Word leaves it
# once the replacement text has been typed. Adding it moves the bookmark
# here from wherever it previously lived.
$postSplit = $d.Range($digitEnd, $digitEnd)
$d.Bookmarks.Add("_GoBack", $postSplit) | Out-Null

# Replace the single "2" character with "3".
$digitRange = $d.Range($digitStart, $digitEnd)
$digitRange.Text = "3"

# Drop the scaffolding bookmark now that the run boundary it enforced is
# baked into the document structure.
$d.Bookmarks("zzTempSplit").Delete()
